# refactor(import): Refactor import programs and inventaries
#
# - Add a "test" value in H5 (new shared string).
# - Clear the (stray) explicit number-format style on E2, returning it to
#   the default/general style.
# - Update the active selection to I14, matching the author's last
#   selection when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = "test"

$ws.Range("E2").ClearFormats() | Out-Null

$ws.Range("I14").Select() | Out-Null
